$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.810.58"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.317.99"
$ws.Range("E3").Value = "  +4.64%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "268.63"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.42"
$ws.Range("E6").Value = "  +7.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.82"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.99"
$ws.Range("E12").Value = "  +6.48%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.664.27"
$ws.Range("E14").Value = "  +4.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.31"
$ws.Range("E15").Value = "  +4.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").Value = "  +9.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.324.45"
$ws.Range("E17").Value = "  +5.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.760.53"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.30"
$ws.Range("E20").Value = "  +5.17%  "
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.72"
$ws.Range("E22").Value = "  +4.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.27"
$ws.Range("E23").Value = "  -4.06%  "
$ws.Range("E24").Value = "  +9.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("E26").Value = "  -6.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.24"
$ws.Range("E27").Value = "  +4.32%  "
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("E29").Value = "  -4.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.94"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.57"
$ws.Range("E31").Value = "  +9.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "172.49"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0893"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.111"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.52"
$ws.Range("E37").Value = "  +3.98%  "
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.237"
$ws.Range("E40").Value = "  +16.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.30"
$ws.Range("E41").Value = "  +9.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.20"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.33"
$ws.Range("E43").Value = "  +18.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.28"
$ws.Range("E45").Value = "  -5.55%  "
$ws.Range("E46").Value = "  +6.68%  "
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.43"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.542.84"
$ws.Range("E50").Value = "  +4.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.432"
$ws.Range("E51").Value = "  -1.32%  "
